$d = $word.ActiveDocument

# The justification paragraph's wording changes from "... the developer
# will use is iterative." to "... the developer will use is agile."
[void]$d.Content.Find.Execute("iterative", $false, $false, $false, $false, $false, `
    $true, 1, $false, "agile", 2)

# The run hosting the WBS-chart inline picture (the first InlineShape,
# originally embedded via r:embed="rId7") picks up <w:noProof/> in its
# run properties.
$d.InlineShapes.Item(1).Range.NoProofing = $true
